$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Range("A33").Value = "desarrollar num correlativa RG-010"
$ws.Range("B33").Value = "no comenzado"
$ws.Range("A29").Select()
